$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = 0.9627400996585722
$ws.Range("E2").Value = 0.9627400996585722

# Row 3
$ws.Range("D3").Value = 0.00001984508027770155
$ws.Range("E3").Value = 0.00001984508027770155

# Row 4
$ws.Range("D4").Value = 0.9999999999999998
$ws.Range("E4").Value = 0.9999999999999998

# Row 5
$ws.Range("D5").Value = 0.9999999999999942
$ws.Range("E5").Value = 0.9999999999999942

# Row 6
$ws.Range("D6").Value = 0.9999999999999956
$ws.Range("E6").Value = 0.9999999999999956

# Row 7
$ws.Range("D7").Value = 0.00000000001501740340712968
$ws.Range("E7").Value = 0.9999999999849826

# Row 9
$ws.Range("D9").Value = 0.9999659986334436
$ws.Range("E9").Value = 0.00003400136655640207

# Row 10
$ws.Range("D10").Value = 0.0000000005333455705050096
$ws.Range("E10").Value = 0.9999999994666544

# Row 11
$ws.Range("D11").Value = 0.00000000002732198292237394
$ws.Range("E11").Value = 0.999999999972678
$ws.Range("F11").Value = 17.5637035369873
$ws.Range("G11").Value = 0.3
